# Swap columns C ("codeforiati:group-name") and D ("codeforiati:group-code")
# for the header row and all data rows, so that column C becomes
# "codeforiati:group-code" and column D becomes "codeforiati:group-name".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$colC = $ws.Range("C1:C$lastRow")
$colD = $ws.Range("D1:D$lastRow")

$valuesC = $colC.Value()
$valuesD = $colD.Value()

$colC.Value = $valuesD
$colD.Value = $valuesC
